$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    3 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 16.84135478251809)
    4 = @(0.0001021024915524027, 0.002571899574220771, 22.3905356188092, 10.19245300693656, 1, 32.58566262781153)
    5 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    6 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    7 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
